$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12 for error code 32: "Failed to open landscape costs
# file in readLandChange()" (shifts existing rows 12-48 down by one).
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = 32
$ws.Range("B12").Value = "Failed to open landscape costs file in readLandChange()"

# Insert a new row at 17 for error code 38: "Invalid SMS cost read from
# costs file in readLandChange()" (shifts rows 17-49 down by one).
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = 38
$ws.Range("B17").Value = "Invalid SMS cost read from costs file in readLandChange()"

# Renumber the two pre-existing rows whose error codes shifted down by one
# (31 -> 30, and 32 -> 31); their messages are unchanged.
$ws.Range("A10").Value = 30
$ws.Range("A11").Value = 31
